# Excel importer bug fixed: column B values were being written as text
# (shared strings) instead of numeric values. Replace the bogus text
# entries in B3:B7 with the correct numeric readings, and leave the
# cursor where the fix was last checked (D8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6575567
$ws.Range("B4").Value = 56765
$ws.Range("B5").Value = 567765576
$ws.Range("B6").Value = 567765
$ws.Range("B7").Value = 567765756

$ws.Range("D8").Select()
